$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 0
$ws.Range("D2").Value = [double]"2.075288033717824E-15"
$ws.Range("E2").Value = [double]"2.075288033717824E-15"

# Row 3 - Control 6
$ws.Range("D3").Value = 0.001415447640223733
$ws.Range("E3").Value = 0.001415447640223733

# Row 4 - Control 9
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = 0.9999973526540379
$ws.Range("E4").Value = 0.9999973526540379

# Row 5 - Control 24
$ws.Range("D5").Value = 0.9999999994247712
$ws.Range("E5").Value = 0.9999999994247712

# Row 6 - Control 32
$ws.Range("D6").Value = 0.9999999993924611
$ws.Range("E6").Value = 0.9999999993924611

# Row 7 - MDD 4
$ws.Range("D7").Value = 0.9979599796227198
$ws.Range("E7").Value = 0.002040020377280238

# Row 8 - MDD 32
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.4314795984705317
$ws.Range("E8").Value = 0.5685204015294683

# Row 10 - MDD 44
$ws.Range("D10").Value = 0.000264710588814839
$ws.Range("E10").Value = 0.9997352894111852

# Row 11 - MDD 31
$ws.Range("D11").Value = 0.99938655323331
$ws.Range("E11").Value = 0.0006134467666899512
$ws.Range("F11").Value = 6.442128658294678
$ws.Range("G11").Value = 0.5
